$d = $word.ActiveDocument
$d.Content.Find.Execute("Unreal Engine 4.20.2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Unreal Engine 4.20.3", 2)
